# Update cryptocurrency price/volume figures (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.249.15"
$ws.Range("E2").Value = "'  +3.19%  "
$ws.Range("D3").Value = "'1.908.91"
$ws.Range("E3").Value = "'  +0.27%  "
$ws.Range("E4").Value = "'  -0.35%  "
$ws.Range("D5").Value = "'326.44"
$ws.Range("E5").Value = "'  +3.49%  "
$ws.Range("E7").Value = "'  +0.52%  "
$ws.Range("D8").Value = "'0.4030"
$ws.Range("E8").Value = "'  +2.21%  "
$ws.Range("D9").Value = "'0.08488"
$ws.Range("E9").Value = "'  +0.25%  "
$ws.Range("D10").Value = "'42.75"
$ws.Range("E10").Value = "'  +0.53%  "
$ws.Range("E11").Value = "'  -0.14%  "
$ws.Range("D12").Value = "'23.50"
$ws.Range("E12").Value = "'  +14.07%  "
$ws.Range("D13").Value = "'6.464"
$ws.Range("E13").Value = "'  +3.16%  "
$ws.Range("D14").Value = "'1.909.15"
$ws.Range("E14").Value = "'  +0.29%  "
$ws.Range("D15").Value = "'7.360"
$ws.Range("E15").Value = "'  -0.02%  "
$ws.Range("E16").Value = "'  -0.32%  "
$ws.Range("D17").Value = "'95.16"
$ws.Range("E17").Value = "'  +1.91%  "
$ws.Range("D18").Value = "'0.00001114"
$ws.Range("E18").Value = "'  +0.60%  "
$ws.Range("D19").Value = "'0.06679"
$ws.Range("E19").Value = "'  -0.81%  "
$ws.Range("D20").Value = "'18.36"
$ws.Range("E20").Value = "'  +2.31%  "
$ws.Range("E21").Value = "'  -0.35%  "
$ws.Range("D22").Value = "'5.995"
$ws.Range("E22").Value = "'  -0.81%  "
$ws.Range("D23").Value = "'30.242.09"
$ws.Range("E23").Value = "'  +3.17%  "
$ws.Range("D24").Value = "'11.30"
$ws.Range("E24").Value = "'  +1.03%  "
$ws.Range("D25").Value = "'2.223"
$ws.Range("E25").Value = "'  +0.05%  "
$ws.Range("D26").Value = "'2.128.15"
$ws.Range("E26").Value = "'  +0.38%  "
$ws.Range("D27").Value = "'21.71"
$ws.Range("E27").Value = "'  +3.74%  "
$ws.Range("D28").Value = "'161.47"
$ws.Range("E28").Value = "'  +0.94%  "
$ws.Range("D29").Value = "'2.401"
$ws.Range("E29").Value = "'  -2.40%  "
$ws.Range("D30").Value = "'129.64"
$ws.Range("E30").Value = "'  +1.27%  "
$ws.Range("D31").Value = "'1.100"
$ws.Range("E31").Value = "'  +3.62%  "
$ws.Range("E32").Value = "'  +0.96%  "
$ws.Range("D33").Value = "'6.051"
$ws.Range("E33").Value = "'  -0.69%  "
$ws.Range("E34").Value = "'  +2.99%  "
$ws.Range("D35").Value = "'0.02507"
$ws.Range("D36").Value = "'0.06589"
$ws.Range("E36").Value = "'  -0.15%  "
$ws.Range("D37").Value = "'0.2216"
$ws.Range("E37").Value = "'  +0.70%  "
$ws.Range("E38").Value = "'  +1.85%  "
$ws.Range("D39").Value = "'1.239"
$ws.Range("E39").Value = "'  -0.08%  "
$ws.Range("D40").Value = "'11.95"
$ws.Range("E40").Value = "'  +5.83%  "
$ws.Range("D41").Value = "'8.833"
$ws.Range("E41").Value = "'  -3.33%  "
$ws.Range("E42").Value = "'  -0.03%  "
$ws.Range("D43").Value = "'1.238"
$ws.Range("E43").Value = "'  +0.20%  "
$ws.Range("D44").Value = "'0.6136"
$ws.Range("E44").Value = "'  +1.26%  "
$ws.Range("D45").Value = "'13.21"
$ws.Range("E45").Value = "'  -0.10%  "
$ws.Range("D46").Value = "'3.718"
$ws.Range("E46").Value = "'  +1.01%  "
$ws.Range("E48").Value = "'  +1.28%  "
$ws.Range("D49").Value = "'125.14"
$ws.Range("E49").Value = "'  +1.58%  "
$ws.Range("D50").Value = "'1.162"
$ws.Range("E50").Value = "'  -1.20%  "
$ws.Range("D51").Value = "'79.49"
$ws.Range("E51").Value = "'  +2.09%  "
